$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93; this pushes existing rows 93-150 down to 94-151.
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new data record.
$ws.Cells.Item(93, 1).Value = 7
$ws.Cells.Item(93, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(93, 3).Value = "Ñuble"
$ws.Cells.Item(93, 4).Value = 44438
$ws.Cells.Item(93, 5).Value = 16
$ws.Cells.Item(93, 6).Value = 100112009
$ws.Cells.Item(93, 7).Value = "Acelga"
$ws.Cells.Item(93, 8).Value = "Sin especificar"
$ws.Cells.Item(93, 9).Value = "Primera"
$ws.Cells.Item(93, 10).Value = 160
$ws.Cells.Item(93, 11).Value = 400
$ws.Cells.Item(93, 12).Value = 450
$ws.Cells.Item(93, 13).Value = 425
$ws.Cells.Item(93, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(93, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(93, 16).Value = 283
$ws.Cells.Item(93, 17).Value = 1.5
$ws.Cells.Item(93, 18).Value = "Hortaliza"
